$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.25517891189561936
$ws.Range("A2").Value = -0.0059999999642599278
$ws.Range("A3").Value = -0.0039999999664210861
$ws.Range("A4").Value = -0.0079999999392317278
$ws.Range("A5").Value = -0.0029999999627481344
$ws.Range("A6").Value = -0.0019999999553270698
$ws.Range("A7").Value = -0.009999999909534818
$ws.Range("A8").Value = -0.02130194861916701
$ws.Range("A9").Value = -0.0019999999538682367
$ws.Range("A10").Value = -0.0019999999544442204
$ws.Range("A11").Value = -0.0029999999488099505
$ws.Range("A12").Value = -0.0034999999463907194
$ws.Range("A13").Value = -0.0034999999486133859
$ws.Range("A14").Value = -0.0079999999242224007
$ws.Range("A15").Value = 0.03669175600679786
$ws.Range("A16").Value = -0.0019999999589486173
$ws.Range("A17").Value = -0.0019999999580253558
$ws.Range("A18").Value = -0.0039999999465223368
$ws.Range("A19").Value = -0.047642278871539823
$ws.Range("A20").Value = -0.0039999999728905777
$ws.Range("A21").Value = -0.0039999999725948143
$ws.Range("A22").Value = -0.0039999999723647761
$ws.Range("A23").Value = -0.0049999999589225297
$ws.Range("A24").Value = -0.019999999870788265
$ws.Range("A25").Value = -0.019999999869092733
$ws.Range("A26").Value = -0.0024999999505528336
$ws.Range("A27").Value = -0.0024999999480139756
$ws.Range("A28").Value = -0.001999999939398478
$ws.Range("A29").Value = -0.0069999999031482574
$ws.Range("A30").Value = -0.059999999600523157
$ws.Range("A31").Value = -0.0069999998961840504
$ws.Range("A32").Value = -0.0099999998788220523
$ws.Range("A33").Value = -0.0039999999124464836
